$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3636.3635
$ws.Range("J43").Value = 2666.6667
$ws.Range("L43").Value = 2666.6667
$ws.Range("N43").Value = -2804.6667

$ws.Range("H53").Value = 519.3
$ws.Range("I53").Value = 307
$ws.Range("J53").Value = 837.75
$ws.Range("K53").Value = 307
$ws.Range("L53").Value = 837.75
$ws.Range("M53").Value = 330
$ws.Range("N53").Value = -2111.75

$ws.Range("H96").Value = 1095
$ws.Range("J96").Value = 316.33334
$ws.Range("L96").Value = 949.0000200000001
$ws.Range("N96").Value = -3695.00002

$ws.Range("H101").Value = 5166.375
$ws.Range("I101").Value = 4020
$ws.Range("K101").Value = 12060
$ws.Range("M101").Value = -10438

$ws.Range("H103").Value = 757.7
$ws.Range("I103").Value = 856.9167
$ws.Range("J103").Value = 608.875
$ws.Range("K103").Value = 2570.7501
$ws.Range("L103").Value = 1826.625
$ws.Range("M103").Value = -1984.7501
$ws.Range("N103").Value = -2998.625

$ws.Range("H106").Value = 2213.6924
$ws.Range("I106").Value = 2323.1667
$ws.Range("J106").Value = 900
$ws.Range("K106").Value = 2323.1667
$ws.Range("L106").Value = 900
$ws.Range("M106").Value = -1692.1667
$ws.Range("N106").Value = -2162

$ws.Range("H107").Value = 466.6
$ws.Range("I107").Value = 466.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 466.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1453.4
$ws.Range("N107").ClearContents()

$ws.Range("H118").Value = 1118.25
$ws.Range("I118").Value = 1118.25
$ws.Range("K118").Value = 3354.75
$ws.Range("M118").Value = -1697.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1472862.2
$ws.Range("I2").Value = 2103046.2
$ws.Range("J2").Value = 2433.3333
$ws.Range("K2").Value = 2103046.2
$ws.Range("L2").Value = 2433.3333
$ws.Range("M2").Value = -2102933.2
$ws.Range("N2").Value = -2659.3333

$ws.Range("H45").Value = 5572.5293
$ws.Range("I45").Value = 6331.154
$ws.Range("K45").Value = 6331.154
$ws.Range("M45").Value = -5954.154

$ws.Range("H61").Value = 142859840
$ws.Range("I61").Value = 166668320
$ws.Range("K61").Value = 166668320
$ws.Range("M61").Value = -166668108

$ws.Range("H88").Value = 2769.7334
$ws.Range("I88").Value = 2641.8333
$ws.Range("J88").Value = 2855
$ws.Range("K88").Value = 2641.8333
$ws.Range("L88").Value = 2855
$ws.Range("M88").Value = -2235.8333
$ws.Range("N88").Value = -3667

$ws.Range("H91").Value = 2769.7334
$ws.Range("I91").Value = 2641.8333
$ws.Range("J91").Value = 2855
$ws.Range("K91").Value = 2641.8333
$ws.Range("L91").Value = 2855
$ws.Range("M91").Value = -1237.8333
$ws.Range("N91").Value = -5663

$ws.Range("H107").Value = 60228
$ws.Range("J107").Value = 60228
$ws.Range("L107").Value = 60228
$ws.Range("N107").Value = -67908

$ws.Range("H116").Value = 1472862.2
$ws.Range("I116").Value = 2103046.2
$ws.Range("J116").Value = 2433.3333
$ws.Range("K116").Value = 2103046.2
$ws.Range("L116").Value = 2433.3333
$ws.Range("M116").Value = -2100752.2
$ws.Range("N116").Value = -7021.3333

$ws.Range("H136").Value = 142859840
$ws.Range("I136").Value = 166668320
$ws.Range("K136").Value = 500004960
$ws.Range("M136").Value = -500002410

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1472862.2
$ws.Range("I3").Value = 2103046.2
$ws.Range("J3").Value = 2433.3333
$ws.Range("K3").Value = 2103046.2
$ws.Range("L3").Value = 2433.3333
$ws.Range("M3").Value = -2102932.2
$ws.Range("N3").Value = -2661.3333

$ws.Range("H4").Value = 3870.1333
$ws.Range("I4").Value = 183.33333
$ws.Range("K4").Value = 183.33333
$ws.Range("M4").Value = -68.33332999999999

$ws.Range("H64").Value = 695.8461
$ws.Range("I64").Value = 501.25
$ws.Range("K64").Value = 501.25
$ws.Range("M64").Value = -276.25

$ws.Range("H67").Value = 695.8461
$ws.Range("I67").Value = 501.25
$ws.Range("K67").Value = 501.25
$ws.Range("M67").Value = 278.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 9774.286
$ws.Range("I70").Value = 2871.2222
$ws.Range("J70").Value = 22199.8
$ws.Range("K70").Value = 8613.6666
$ws.Range("L70").Value = 66599.39999999999
$ws.Range("M70").Value = -8298.6666
$ws.Range("N70").Value = -67229.39999999999

$ws.Range("H73").Value = 9774.286
$ws.Range("I73").Value = 2871.2222
$ws.Range("J73").Value = 22199.8
$ws.Range("K73").Value = 8613.6666
$ws.Range("L73").Value = 66599.39999999999
$ws.Range("M73").Value = -7521.6666
$ws.Range("N73").Value = -68783.39999999999

$ws.Range("H99").Value = 5746.25
$ws.Range("I99").Value = 1278
$ws.Range("J99").Value = 10214.5
$ws.Range("K99").Value = 3834
$ws.Range("L99").Value = 30643.5
$ws.Range("M99").Value = -1588
$ws.Range("N99").Value = -35135.5

$ws.Range("H140").Value = 451.6
$ws.Range("I140").Value = 451.6
$ws.Range("K140").Value = 1354.8
$ws.Range("M140").Value = 3825.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1708.5
$ws.Range("I43").Value = 1708.5
$ws.Range("K43").Value = 1708.5
$ws.Range("M43").Value = -1557.5

$ws.Range("H97").Value = 1464.697
$ws.Range("I97").Value = 1268.5
$ws.Range("J97").Value = 1987.8889
$ws.Range("K97").Value = 1268.5
$ws.Range("L97").Value = 1987.8889
$ws.Range("M97").Value = -772.5
$ws.Range("N97").Value = -2979.8889

$ws.Range("H124").Value = 49998.5
$ws.Range("J124").Value = 49998.5
$ws.Range("L124").Value = 49998.5
$ws.Range("N124").Value = -59818.5

$ws.Range("H126").Value = 3969
$ws.Range("I126").Value = 3969
$ws.Range("K126").Value = 11907
$ws.Range("M126").Value = -9437

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2905.6667
$ws.Range("I22").Value = 2970.4167
$ws.Range("J22").Value = 2646.6667
$ws.Range("K22").Value = 2970.4167
$ws.Range("L22").Value = 2646.6667
$ws.Range("M22").Value = -2675.4167
$ws.Range("N22").Value = -3236.6667

$ws.Range("H27").Value = 2905.6667
$ws.Range("I27").Value = 2970.4167
$ws.Range("J27").Value = 2646.6667
$ws.Range("K27").Value = 2970.4167
$ws.Range("L27").Value = 2646.6667
$ws.Range("M27").Value = -2863.4167
$ws.Range("N27").Value = -2860.6667

$ws.Range("H40").Value = 2952.8
$ws.Range("J40").Value = 2999.5
$ws.Range("L40").Value = 2999.5
$ws.Range("N40").Value = -3271.5

$ws.Range("H46").Value = 1707.95
$ws.Range("I46").Value = 1821.1177
$ws.Range("K46").Value = 1821.1177
$ws.Range("M46").Value = -1633.1177

$ws.Range("H55").Value = 484.8889
$ws.Range("I55").Value = 174.3
$ws.Range("J55").Value = 873.125
$ws.Range("K55").Value = 174.3
$ws.Range("L55").Value = 873.125
$ws.Range("M55").Value = -1.300000000000011
$ws.Range("N55").Value = -1219.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3434.12
$ws.Range("I96").Value = 1353.6
$ws.Range("J96").Value = 3954.25
$ws.Range("K96").Value = 1353.6
$ws.Range("L96").Value = 3954.25
$ws.Range("M96").Value = 19.40000000000009
$ws.Range("N96").Value = -6700.25

$ws.Range("H103").Value = 26200
$ws.Range("J103").Value = 26200
$ws.Range("L103").Value = 26200
$ws.Range("N103").Value = -28544

$ws.Range("H132").Value = 14290930
$ws.Range("I132").Value = 20002250
$ws.Range("K132").Value = 60006750
$ws.Range("M132").Value = -60004220

$ws.Range("H136").Value = 17243062
$ws.Range("I136").Value = 20834896
$ws.Range("K136").Value = 62504688
$ws.Range("M136").Value = -62502138
